# Auto-generated script applying scheduled market-data refresh to Typhon_Profits workbook.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) per-sheet per-row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5469.85
$ws.Range("I62").Value = 4273.1816
$ws.Range("J62").Value = 6932.4443
$ws.Range("K62").Value = 4273.1816
$ws.Range("L62").Value = 6932.4443
$ws.Range("M62").Value = -3649.1816
$ws.Range("N62").Value = -8180.4443
$ws.Range("H65").Value = 5469.85
$ws.Range("I65").Value = 4273.1816
$ws.Range("J65").Value = 6932.4443
$ws.Range("K65").Value = 21365.908
$ws.Range("L65").Value = 34662.2215
$ws.Range("M65").Value = -18245.908
$ws.Range("N65").Value = -40902.2215
$ws.Range("H137").Value = 1829.1666
$ws.Range("J137").Value = 1535.5714
$ws.Range("L137").Value = 4606.7142
$ws.Range("N137").Value = -9706.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5220.653
$ws.Range("I32").Value = 4451.892
$ws.Range("J32").Value = 12359.143
$ws.Range("K32").Value = 4451.892
$ws.Range("L32").Value = 12359.143
$ws.Range("M32").Value = -4164.892
$ws.Range("N32").Value = -12933.143
$ws.Range("H63").Value = 1954891
$ws.Range("J63").Value = 7813626.5
$ws.Range("L63").Value = 7813626.5
$ws.Range("N63").Value = -7814998.5
$ws.Range("H66").Value = 1954891
$ws.Range("J66").Value = 7813626.5
$ws.Range("L66").Value = 39068132.5
$ws.Range("N66").Value = -39074996.5
$ws.Range("H74").Value = 50001950
$ws.Range("I74").Value = 125000664
$ws.Range("K74").Value = 125000664
$ws.Range("M74").Value = -124999790
$ws.Range("H77").Value = 50001950
$ws.Range("I77").Value = 125000664
$ws.Range("K77").Value = 625003320
$ws.Range("M77").Value = -624998952
$ws.Range("H88").Value = 112672.336
$ws.Range("I88").Value = 1647.6666
$ws.Range("J88").Value = 168184.67
$ws.Range("K88").Value = 1647.6666
$ws.Range("L88").Value = 168184.67
$ws.Range("M88").Value = -1241.6666
$ws.Range("N88").Value = -168996.67
$ws.Range("H91").Value = 112672.336
$ws.Range("I91").Value = 1647.6666
$ws.Range("J91").Value = 168184.67
$ws.Range("K91").Value = 1647.6666
$ws.Range("L91").Value = 168184.67
$ws.Range("M91").Value = -243.6666
$ws.Range("N91").Value = -170992.67

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1583.8334
$ws.Range("I86").Value = 1454.375
$ws.Range("J86").Value = 1713.2916
$ws.Range("K86").Value = 1454.375
$ws.Range("L86").Value = 1713.2916
$ws.Range("M86").Value = -331.375
$ws.Range("N86").Value = -3959.2916
$ws.Range("H89").Value = 1583.8334
$ws.Range("I89").Value = 1454.375
$ws.Range("J89").Value = 1713.2916
$ws.Range("K89").Value = 7271.875
$ws.Range("L89").Value = 8566.458000000001
$ws.Range("M89").Value = -1655.875
$ws.Range("N89").Value = -19798.458
$ws.Range("H94").Value = 1005.9
$ws.Range("I94").Value = 794.1429000000001
$ws.Range("K94").Value = 794.1429000000001
$ws.Range("M94").Value = -343.1429000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3448.2927
$ws.Range("I31").Value = 3992.7778
$ws.Range("J31").Value = 3295.1562
$ws.Range("K31").Value = 3992.7778
$ws.Range("L31").Value = 3295.1562
$ws.Range("M31").Value = -3697.7778
$ws.Range("N31").Value = -3885.1562
$ws.Range("H34").Value = 3448.2927
$ws.Range("I34").Value = 3992.7778
$ws.Range("J34").Value = 3295.1562
$ws.Range("K34").Value = 3992.7778
$ws.Range("L34").Value = 3295.1562
$ws.Range("M34").Value = -3790.7778
$ws.Range("N34").Value = -3699.1562
$ws.Range("H58").Value = 18280.232
$ws.Range("I58").Value = 1566.2667
$ws.Range("K58").Value = 1566.2667
$ws.Range("M58").Value = -1363.2667
$ws.Range("H62").Value = 55559136
$ws.Range("J62").Value = 3475
$ws.Range("L62").Value = 3475
$ws.Range("N62").Value = -4723
$ws.Range("H65").Value = 55559136
$ws.Range("J65").Value = 3475
$ws.Range("L65").Value = 17375
$ws.Range("N65").Value = -23615
$ws.Range("H132").Value = 2498.6765
$ws.Range("I132").Value = 1877.174
$ws.Range("J132").Value = 3798.182
$ws.Range("K132").Value = 5631.522
$ws.Range("L132").Value = 11394.546
$ws.Range("M132").Value = -3101.522
$ws.Range("N132").Value = -16454.546
$ws.Range("H134").Value = 1139.742
$ws.Range("I134").Value = 1034.5186
$ws.Range("K134").Value = 3103.5558
$ws.Range("M134").Value = -568.5558000000001
$ws.Range("H136").Value = 18280.232
$ws.Range("I136").Value = 1566.2667
$ws.Range("K136").Value = 4698.800099999999
$ws.Range("M136").Value = -2148.800099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1853.3636
$ws.Range("I68").Value = 1366.6666
$ws.Range("J68").Value = 2035.875
$ws.Range("K68").Value = 4099.9998
$ws.Range("L68").Value = 6107.625
$ws.Range("M68").Value = -3288.9998
$ws.Range("N68").Value = -7729.625
$ws.Range("H71").Value = 1853.3636
$ws.Range("I71").Value = 1366.6666
$ws.Range("J71").Value = 2035.875
$ws.Range("K71").Value = 12299.9994
$ws.Range("L71").Value = 18322.875
$ws.Range("M71").Value = -8243.999400000001
$ws.Range("N71").Value = -26434.875
$ws.Range("H81").Value = 4334
$ws.Range("I81").Value = 3000
$ws.Range("J81").Value = 4539.231
$ws.Range("K81").Value = 9000
$ws.Range("L81").Value = 13617.693
$ws.Range("M81").Value = -7877
$ws.Range("N81").Value = -15863.693
$ws.Range("H84").Value = 4334
$ws.Range("I84").Value = 3000
$ws.Range("J84").Value = 4539.231
$ws.Range("K84").Value = 27000
$ws.Range("L84").Value = 40853.079
$ws.Range("M84").Value = -21384
$ws.Range("N84").Value = -52085.079
$ws.Range("H98").Value = 728
$ws.Range("I98").Value = 799.3333
$ws.Range("K98").Value = 2397.9999
$ws.Range("M98").Value = -899.9998999999998
$ws.Range("H113").Value = 736.625
$ws.Range("I113").Value = 601.1429000000001
$ws.Range("J113").Value = 842
$ws.Range("K113").Value = 1803.4287
$ws.Range("L113").Value = 2526
$ws.Range("M113").Value = 366.5712999999998
$ws.Range("N113").Value = -6866
$ws.Range("H121").Value = 1620
$ws.Range("J121").Value = 1733.3334
$ws.Range("L121").Value = 5200.0002
$ws.Range("N121").Value = -7820.0002
$ws.Range("H131").Value = 114409.02
$ws.Range("I131").Value = 903.3333
$ws.Range("J131").Value = 118415.11
$ws.Range("K131").Value = 2709.9999
$ws.Range("L131").Value = 355245.33
$ws.Range("M131").Value = 2330.0001
$ws.Range("N131").Value = -365325.33

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3060.8928
$ws.Range("I80").Value = 2099.7693
$ws.Range("J80").Value = 3893.8667
$ws.Range("K80").Value = 2099.7693
$ws.Range("L80").Value = 3893.8667
$ws.Range("M80").Value = -1101.7693
$ws.Range("N80").Value = -5889.8667
$ws.Range("H83").Value = 3060.8928
$ws.Range("I83").Value = 2099.7693
$ws.Range("J83").Value = 3893.8667
$ws.Range("K83").Value = 10498.8465
$ws.Range("L83").Value = 19469.3335
$ws.Range("M83").Value = -5506.8465
$ws.Range("N83").Value = -29453.3335
$ws.Range("H122").Value = 41667640
$ws.Range("I122").Value = 15152463
$ws.Range("J122").Value = 100001020
$ws.Range("K122").Value = 45457389
$ws.Range("L122").Value = 300003060
$ws.Range("M122").Value = -45454939
$ws.Range("N122").Value = -300007960
$ws.Range("H126").Value = 5521.3335
$ws.Range("I126").Value = 4400
$ws.Range("J126").Value = 7764
$ws.Range("K126").Value = 13200
$ws.Range("L126").Value = 23292
$ws.Range("M126").Value = -10730
$ws.Range("N126").Value = -28232

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2945.4517
$ws.Range("I40").Value = 2441.5908
$ws.Range("J40").Value = 4177.1113
$ws.Range("K40").Value = 2441.5908
$ws.Range("L40").Value = 4177.1113
$ws.Range("M40").Value = -2305.5908
$ws.Range("N40").Value = -4449.1113
$ws.Range("H68").Value = 2149.25
$ws.Range("I68").Value = 1716.3334
$ws.Range("J68").Value = 2582.1667
$ws.Range("K68").Value = 1716.3334
$ws.Range("L68").Value = 2582.1667
$ws.Range("M68").Value = -967.3334
$ws.Range("N68").Value = -4080.1667
$ws.Range("H71").Value = 2149.25
$ws.Range("I71").Value = 1716.3334
$ws.Range("J71").Value = 2582.1667
$ws.Range("K71").Value = 8581.666999999999
$ws.Range("L71").Value = 12910.8335
$ws.Range("M71").Value = -4837.666999999999
$ws.Range("N71").Value = -20398.8335
$ws.Range("H122").Value = 855540.4
$ws.Range("I122").Value = 1309656.2
$ws.Range("K122").Value = 3928968.6
$ws.Range("M122").Value = -3926518.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 997.5357
$ws.Range("I126").Value = 1053.6086
$ws.Range("K126").Value = 3160.8258
$ws.Range("M126").Value = -690.8258000000001
$ws.Range("H136").Value = 21068720
$ws.Range("I136").Value = 29494114
$ws.Range("K136").Value = 88482342
$ws.Range("M136").Value = -88479792

